$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").EntireColumn.Delete()

$ws.Range("A2").Value = "0000_ORG_Phone_Scenario"
$ws.Range("B2").Value = "Verify the ORG_Phone Info"
$ws.Range("C2").Value = "Verify_All_Buckets_ORG_PHONE"
$ws.Range("D2").Value = "100"

$ws.Range("H19").Select()
